$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "42.258.59"
Set-TextValue $ws.Range("E2") "  +0.53%  "

Set-TextValue $ws.Range("D3") "2.284.25"
Set-TextValue $ws.Range("E3") "  -0.78%  "

Set-TextValue $ws.Range("E4") "  -0.06%  "

Set-TextValue $ws.Range("D5") "325.33"
Set-TextValue $ws.Range("E5") "  +2.62%  "

Set-TextValue $ws.Range("D6") "102.88"
Set-TextValue $ws.Range("E6") "  -1.86%  "

Set-TextValue $ws.Range("E7") "  -0.88%  "

Set-TextValue $ws.Range("E8") "  +0.09%  "

Set-TextValue $ws.Range("D9") "0.607"
Set-TextValue $ws.Range("E9") "  -0.33%  "

Set-TextValue $ws.Range("D10") "39.88"
Set-TextValue $ws.Range("E10") "  +0.44%  "

Set-TextValue $ws.Range("D11") "0.0903"
Set-TextValue $ws.Range("E11") "  -0.71%  "

Set-TextValue $ws.Range("D12") "8.35"
Set-TextValue $ws.Range("E12") "  -1.07%  "

Set-TextValue $ws.Range("E13") "  +0.00%  "

Set-TextValue $ws.Range("E14") "  -0.83%  "

Set-TextValue $ws.Range("D15") "15.11"
Set-TextValue $ws.Range("E15") "  -2.16%  "

Set-TextValue $ws.Range("D16") "2.631.74"
Set-TextValue $ws.Range("E16") "  -0.75%  "

Set-TextValue $ws.Range("D17") "2.285.92"
Set-TextValue $ws.Range("E17") "  -0.60%  "

Set-TextValue $ws.Range("D18") "42.250.40"
Set-TextValue $ws.Range("E18") "  +0.28%  "

Set-TextValue $ws.Range("E19") "  -5.03%  "

Set-TextValue $ws.Range("E20") "  -0.52%  "

Set-TextValue $ws.Range("D21") "12.96"
Set-TextValue $ws.Range("E21") "  +29.98%  "

Set-TextValue $ws.Range("E22") "  +2.57%  "

Set-TextValue $ws.Range("D23") "73.04"
Set-TextValue $ws.Range("E23") "  -0.88%  "

Set-TextValue $ws.Range("D24") "267.72"
Set-TextValue $ws.Range("E24") "  -6.88%  "

Set-TextValue $ws.Range("E25") "  -2.95%  "

Set-TextValue $ws.Range("E26") "  -0.03%  "

Set-TextValue $ws.Range("E27") "  -0.99%  "

Set-TextValue $ws.Range("E28") "  +4.37%  "

Set-TextValue $ws.Range("B29") "EthereumClassic"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D29") "22.43"
Set-TextValue $ws.Range("E29") "  -4.26%  "

Set-TextValue $ws.Range("B30") "InjectiveProtocol"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D30") "37.79"
Set-TextValue $ws.Range("E30") "  +6.45%  "

Set-TextValue $ws.Range("D31") "164.19"
Set-TextValue $ws.Range("E31") "  -0.53%  "

Set-TextValue $ws.Range("D32") "6.16"
Set-TextValue $ws.Range("E32") "  +3.98%  "

Set-TextValue $ws.Range("D33") "0.0876"
Set-TextValue $ws.Range("E33") "  -0.80%  "

Set-TextValue $ws.Range("E34") "  +0.56%  "

Set-TextValue $ws.Range("B35") "Kaspa"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D35") "0.113"
Set-TextValue $ws.Range("E35") "  -2.70%  "

Set-TextValue $ws.Range("B36") "WEMIXToken"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D36") "2.50"
Set-TextValue $ws.Range("E36") "  -14.52%  "

Set-TextValue $ws.Range("E37") "  -0.92%  "

Set-TextValue $ws.Range("D38") "0.0353"
Set-TextValue $ws.Range("E38") "  +0.07%  "

Set-TextValue $ws.Range("D39") "3.74"
Set-TextValue $ws.Range("E39") "  +2.95%  "

Set-TextValue $ws.Range("E40") "  -6.80%  "

Set-TextValue $ws.Range("E41") "  +1.44%  "

Set-TextValue $ws.Range("D42") "69.50"
Set-TextValue $ws.Range("E42") "  -2.51%  "

Set-TextValue $ws.Range("E43") "  -0.37%  "

Set-TextValue $ws.Range("E44") "  -1.16%  "

Set-TextValue $ws.Range("D45") "90.95"
Set-TextValue $ws.Range("E45") "  -11.08%  "

Set-TextValue $ws.Range("D46") "12.28"
Set-TextValue $ws.Range("E46") "  +1.28%  "

Set-TextValue $ws.Range("B47") "ordi"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws.Range("D47") "80.13"
Set-TextValue $ws.Range("E47") "  +2.09%  "

Set-TextValue $ws.Range("B48") "Aave"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D48") "112.70"
Set-TextValue $ws.Range("E48") "  -2.98%  "

Set-TextValue $ws.Range("E49") "  -2.92%  "

Set-TextValue $ws.Range("E50") "  -2.51%  "

Set-TextValue $ws.Range("B51") "MinaProtocolToken"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
Set-TextValue $ws.Range("D51") "1.52"
Set-TextValue $ws.Range("E51") "  +8.90%  "

